$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "Player Info" worksheet before the existing "ODI Batting"
#    sheet so the final order is: Player Info, ODI Batting, ODI Bowling.
#    NOTE: worksheet references returned by this host are positional, so we
#    re-fetch sheets by name after any structural change (add/move/rename).
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($beforeSheet)
$playerInfo.Name = "Player Info"

# Re-fetch a fresh reference to the new sheet by its final name.
$playerInfo = $wb.Worksheets.Item("Player Info")

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold / centered / thin-bordered look used by the other sheet headers
$headerRng = $playerInfo.Range("A1:D1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

# Data row
$playerInfo.Range("A2").Value = "'6039"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Saud Shakeel"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

# ---------------------------------------------------------------------------
# 2. "ODI Batting" sheet: rename MATCH_CARD_LINK -> MATCH_CODE and replace the
#    full scorecard URLs with just the bare match codes.
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{ 2 = "4472"; 3 = "4473"; 4 = "4476"; 5 = "4564"; 6 = "4565" }
foreach ($row in $battingCodes.Keys) {
    $cell = $battingSheet.Range("D$row")
    $cell.Value = "'" + $battingCodes[$row]
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 3. "ODI Bowling" sheet: rename MATCH_CARD_LINK -> MATCH_CODE and replace the
#    full scorecard URLs with just the bare match codes.
# ---------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{ 2 = "4472"; 3 = "4473"; 4 = "4476" }
foreach ($row in $bowlingCodes.Keys) {
    $cell = $bowlingSheet.Range("B$row")
    $cell.Value = "'" + $bowlingCodes[$row]
    $cell.Style = "Normal"
}
